$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column cells keep their original text formatting,
# so numeric-looking strings (e.g. trailing zeros) are not coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.826.74'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.659.94'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.97'
$ws.Range("E5").Value = '  -1.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '327.97'
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.526'
$ws.Range("E7").Value = '  -0.67%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.552'
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.87'
$ws.Range("E10").Value = '  -3.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.98'
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0819'
$ws.Range("E12").Value = '  -0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.59'
$ws.Range("E14").Value = '  +2.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.074.98'
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.653.79'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.868'
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.762.31'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.62'
$ws.Range("E19").Value = '  +2.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.94'
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.72'
$ws.Range("E21").Value = '  -0.89%  '
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '270.23'
$ws.Range("E23").Value = '  -2.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.36'
$ws.Range("E24").Value = '  -4.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.56'
$ws.Range("E25").Value = '  -1.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.24'
$ws.Range("E26").Value = '  -2.74%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.24'
$ws.Range("E28").Value = '  +1.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  -2.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.07'
$ws.Range("E31").Value = '  -5.04%  '
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.51'
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0822'
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.19'
$ws.Range("E35").Value = '  -2.63%  '
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("E38").Value = '  -1.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.17'
$ws.Range("E39").Value = '  +1.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.83'
$ws.Range("E40").Value = '  +7.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '128.15'
$ws.Range("E41").Value = '  +3.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0348'
$ws.Range("E42").Value = '  +8.99%  '
$ws.Range("E43").Value = '  +3.88%  '
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("E45").Value = '  +0.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.065.46'
$ws.Range("E46").Value = '  -1.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.13'
$ws.Range("E47").Value = '  +7.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.20'
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("E49").Value = '  -1.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.32'
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("E51").Value = '  -1.75%  '
